$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data table (header + 17 player rows) in final order.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Spencer Dinwiddie", "PG,SG", "Dallas Mavericks"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat")
)

# Clear out the previous data rows (A2:C19) before writing the new,
# shorter table (rows 2-18). The header row (row 1) keeps its existing
# formatting untouched.
$ws.Range("A2:C19").Clear()

for ($i = 1; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
